$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F - copy E1's formatting (bold/border/alignment style)
# onto F1, then set its text, matching the existing header row styling.
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Per-row timestamps (stored as text, matching the source inlineStr values)
$timestamps = @(
    "2021-10-05 13:39:12.893889",
    "2021-10-05 13:39:12.893902",
    "2021-10-05 13:39:12.893906",
    "2021-10-05 13:39:12.893909",
    "2021-10-05 13:39:12.893913",
    "2021-10-05 13:39:12.893916",
    "2021-10-05 13:39:12.893919",
    "2021-10-05 13:39:12.893922",
    "2021-10-05 13:39:12.893925",
    "2021-10-05 13:39:12.893928",
    "2021-10-05 13:39:12.893931",
    "2021-10-05 13:39:12.893934",
    "2021-10-05 13:39:12.893937",
    "2021-10-05 13:39:12.893940",
    "2021-10-05 13:39:12.893943",
    "2021-10-05 13:39:12.893946",
    "2021-10-05 13:39:12.893949",
    "2021-10-05 13:39:12.893952",
    "2021-10-05 13:39:12.893955",
    "2021-10-05 13:39:12.893958",
    "2021-10-05 13:39:12.893961"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
